$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" '92.602.35'
$ws.Range("E2").Value = '  -4.88%  '

Set-TextValue "D3" '3.326.78'
$ws.Range("E3").Value = '  -3.85%  '

Set-TextValue "D4" '0.999'
$ws.Range("E4").Value = '  -0.19%  '

Set-TextValue "D5" '228.57'
$ws.Range("E5").Value = '  -7.98%  '

Set-TextValue "D6" '622.76'
$ws.Range("E6").Value = '  -4.82%  '

Set-TextValue "D7" '1.34'
$ws.Range("E7").Value = '  -7.28%  '

Set-TextValue "D8" '0.375'
$ws.Range("E8").Value = '  -9.94%  '

Set-TextValue "D9" '0.999'
$ws.Range("E9").Value = '  -0.02%  '

Set-TextValue "D10" '0.907'
$ws.Range("E10").Value = '  -11.02%  '

Set-TextValue "D11" '3.325.52'
$ws.Range("E11").Value = '  -3.73%  '

Set-TextValue "D12" '41.81'
$ws.Range("E12").Value = '  -5.16%  '

Set-TextValue "D13" '0.190'
$ws.Range("E13").Value = '  -8.94%  '

$ws.Range("B14").Value = 'Toncoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue "D14" '5.92'
$ws.Range("E14").Value = '  -3.18%  '

$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue "D15" '92.221.74'
$ws.Range("E15").Value = '  -5.18%  '

Set-TextValue "D16" '3.954.05'
$ws.Range("E16").Value = '  -4.14%  '

Set-TextValue "D17" '0.0000239'
$ws.Range("E17").Value = '  -6.22%  '

Set-TextValue "D18" '7.86'
$ws.Range("E18").Value = '  -9.86%  '

Set-TextValue "D19" '3.322.48'
$ws.Range("E19").Value = '  -4.67%  '

Set-TextValue "D20" '16.76'
$ws.Range("E20").Value = '  -9.33%  '

Set-TextValue "D21" '10.83'
$ws.Range("E21").Value = '  -7.67%  '

Set-TextValue "D22" '485.59'
$ws.Range("E22").Value = '  -5.51%  '

Set-TextValue "D23" '3.26'
$ws.Range("E23").Value = '  -2.49%  '

Set-TextValue "D24" '0.435'
$ws.Range("E24").Value = '  -14.53%  '

Set-TextValue "D25" '0.0000179'
$ws.Range("E25").Value = '  -9.37%  '

Set-TextValue "D26" '6.12'
$ws.Range("E26").Value = '  -9.11%  '

Set-TextValue "D27" '88.98'
$ws.Range("E27").Value = '  -4.36%  '

Set-TextValue "D28" '3.504.12'
$ws.Range("E28").Value = '  -4.15%  '

Set-TextValue "D29" '11.45'
$ws.Range("E29").Value = '  -7.65%  '

$ws.Range("E30").Value = '  +0.05%  '

Set-TextValue "D31" '11.01'
$ws.Range("E31").Value = '  -8.24%  '

$ws.Range("E32").Value = '  -6.06%  '

Set-TextValue "D33" '2.59'
$ws.Range("E33").Value = '  -7.61%  '

Set-TextValue "D34" '1.01'
$ws.Range("E34").Value = '  +1.32%  '

Set-TextValue "D35" '0.168'
$ws.Range("E35").Value = '  -9.90%  '

Set-TextValue "D36" '27.98'
$ws.Range("E36").Value = '  -8.44%  '

Set-TextValue "D37" '0.517'
$ws.Range("E37").Value = '  -10.91%  '

Set-TextValue "D38" '519.93'
$ws.Range("E38").Value = '  +0.86%  '

$ws.Range("E39").Value = '  -0.02%  '

Set-TextValue "D40" '7.23'
$ws.Range("E40").Value = '  -7.15%  '

$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue "D41" '1.34'
$ws.Range("E41").Value = '  -9.08%  '

$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D42" '0.144'
$ws.Range("E42").Value = '  -6.63%  '

Set-TextValue "D43" '0.870'
$ws.Range("E43").Value = '  -3.14%  '

Set-TextValue "D44" '23.97'
$ws.Range("E44").Value = '  -1.43%  '

Set-TextValue "D45" '3.57'
$ws.Range("E45").Value = '  -1.66%  '

Set-TextValue "D46" '1.65'
$ws.Range("E46").Value = '  -3.37%  '

Set-TextValue "D47" '5.29'
$ws.Range("E47").Value = '  -5.56%  '

$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D48" '0.0388'
$ws.Range("E48").Value = '  -8.31%  '

$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D49" '2.09'
$ws.Range("E49").Value = '  -5.16%  '

$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D50" '51.83'
$ws.Range("E50").Value = '  -4.07%  '

Set-TextValue "D51" '7.84'
$ws.Range("E51").Value = '  -7.67%  '
